# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") -- i.e. theme1.xml becomes the
# "Integral" theme and theme2.xml becomes the "Office Theme". Both theme
# parts already share an identical <a:fontScheme> and <a:fmtScheme>; the
# only real content difference between them is the 12-colour
# <a:clrScheme>. The presentation's SlideMaster (and therefore every
# slide) is wired to ppt/theme/theme2.xml, which is the theme part the
# PowerPoint object model exposes here, so we re-point its colour scheme
# at the "Office Theme" palette that theme2.xml is supposed to end up
# with (values are the packed BGR integers PowerPoint's RGB() produces:
# R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # Dark 1       000000
$cs.Item(2).RGB  = 16777215   # Light 1      FFFFFF
$cs.Item(3).RGB  = 6968388    # Dark 2       44546A
$cs.Item(4).RGB  = 15132391   # Light 2      E7E6E6
$cs.Item(5).RGB  = 13998939   # Accent 1     5B9BD5
$cs.Item(6).RGB  = 3243501    # Accent 2     ED7D31
$cs.Item(7).RGB  = 10855845   # Accent 3     A5A5A5
$cs.Item(8).RGB  = 49407      # Accent 4     FFC000
$cs.Item(9).RGB  = 12874308   # Accent 5     4472C4
$cs.Item(10).RGB = 4697456    # Accent 6     70AD47
$cs.Item(11).RGB = 12673797   # Hyperlink    0563C1
$cs.Item(12).RGB = 7491477    # Followed hyperlink  954F72
